$wb = $excel.ActiveWorkbook

$wsDatos  = $wb.Worksheets.Item("Datos")
$wsDatos2 = $wb.Worksheets.Item("Datos2")

# --- Content edits on "Datos" sheet ---------------------------------------
# G3 must become the brand-new value "chipote25" first, so that when it is
# appended to the shared-string table it lands before "Hola".
$wsDatos.Range("G3").Value = "chipote25"

# E2 / E3 change from "¡Hola!" to "Hola". A leading apostrophe is used so the
# COM layer keeps treating the cell with its original (quote-prefixed) cell
# style instead of re-deriving a fresh style for the literal text.
$wsDatos.Range("E2").Value = "'Hola"
$wsDatos.Range("E3").Value = "'Hola"

# --- View / selection state -------------------------------------------------
# "Datos" becomes the active sheet/tab (was "Datos2").
$wsDatos.Activate()
$wsDatos.Range("E10").Select()

# "Datos2" keeps its own selection, just loses the tabSelected flag, which
# Activate() above already took care of.
$wsDatos2.Range("D15").Select()

# restore Datos as the active sheet/selection (selecting on Datos2 above
# switches the active sheet again)
$wsDatos.Activate()
$wsDatos.Range("E10").Select()
